# Apply updated vm_pu values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..N (excluding the empty H column) map to spreadsheet columns:
$cols = @(2, 3, 4, 5, 6, 7, 9, 10, 11, 12, 13, 14)   # B, C, D, E, F, G, I, J, K, L, M, N

# Data rows correspond to worksheet rows 2..25 (row 1 is the header row)
$data = @(
    @(1.02, 1.028772881466093, 1.036835762705145, 1.028699403059596, 1.045671437307821, 1, 1.034203299525416, 1.033923342961585, 1.039628367451097, 1.03151546011113, 1.048439022621327, 1.015328318491058),
    @(1.02, 1.029630425090795, 1.037601584458905, 1.029424753498311, 1.046539170838426, 1, 1.034328516200894, 1.034422108262103, 1.040203989008678, 1.032049023799471, 1.049118097408869, 1.015493305005049),
    @(1.02, 1.030185572535411, 1.038097651324499, 1.029894726611969, 1.047101193420632, 1, 1.034407630309353, 1.034744460519251, 1.040576322825897, 1.032394229605615, 1.049557404743089, 1.015599920067272),
    @(1.02, 1.030419017112085, 1.038306322718284, 1.03009245069078, 1.047337595365715, 1, 1.034440431855096, 1.034879884618082, 1.040732818988229, 1.03253934186963, 1.049742064110021, 1.015644706538097),
    @(1.02, 1.030458216993015, 1.038341366835563, 1.030125658025586, 1.047377295723976, 1, 1.034445912493895, 1.0349026174416, 1.04075909340759, 1.032563706099323, 1.049773067691308, 1.015652224347393),
    @(1.02, 1.03018869159708, 1.038100439113386, 1.02989736803387, 1.047104351735753, 1, 1.034408070406076, 1.034746270428905, 1.040578414066632, 1.032396168652649, 1.049559872272588, 1.015600518642419),
    @(1.02, 1.029062637739656, 1.037094465169329, 1.028944408623588, 1.045964578792618, 1, 1.034246011913771, 1.0340919814081, 1.039822927896171, 1.0316957890878, 1.048668538744202, 1.015384105640024),
    @(1.02, 1.027080439735431, 1.035325939233844, 1.027270019348355, 1.04396037811709, 1, 1.033945859395591, 1.032936169325651, 1.038490717453439, 1.030461340934706, 1.047097204926372, 1.015001688545435),
    @(1.02, 1.025760442870182, 1.034149801503407, 1.02615712131582, 1.042627194345744, 1, 1.033736004294918, 1.032163772787347, 1.037602023669715, 1.029638259929203, 1.046049280647255, 1.014746050945808),
    @(1.02, 1.025189236264238, 1.033641225573495, 1.02567604301224, 1.042050633616547, 1, 1.03364283223122, 1.031828892238135, 1.037217095581108, 1.029281845937305, 1.045595450451307, 1.014635198204602),
    @(1.02, 1.024977120335434, 1.033452424653231, 1.025497473118113, 1.041836582780851, 1, 1.033607878755033, 1.03170443966815, 1.037074099832308, 1.029149457055865, 1.045426868467353, 1.014593998962467),
    @(1.02, 1.025022617381398, 1.033492918229666, 1.025535771327241, 1.04188249242783, 1, 1.033615392002981, 1.031731137998939, 1.037104773622101, 1.029177854937899, 1.045463030268023, 1.014602837400653),
    @(1.02, 1.025171701552084, 1.033625617041317, 1.02566127982938, 1.042032937863597, 1, 1.033639950003587, 1.031818606228373, 1.037205275832803, 1.029270902640898, 1.045581515591319, 1.014631793141279),
    @(1.02, 1.025263564645591, 1.033707391365925, 1.025738626243741, 1.042125646824893, 1, 1.033655035277872, 1.031872489977173, 1.037267196440818, 1.029328232327732, 1.045654517132117, 1.014649630614739),
    @(1.02, 1.025798359694837, 1.034183568898341, 1.026189066193967, 1.04266547403705, 1, 1.033742139359299, 1.032185988800337, 1.037627567715124, 1.029661913759686, 1.046079398479739, 1.014753404543944),
    @(1.02, 1.02613392009924, 1.034482450979839, 1.026471834398765, 1.043004286539321, 1, 1.03379616135169, 1.032382524686278, 1.037853588529709, 1.029871220362048, 1.046345897149979, 1.014818456666856),
    @(1.02, 1.026329681516386, 1.034656851213718, 1.026636846721317, 1.043201979163202, 1, 1.033827449163209, 1.032497119388676, 1.037985411186934, 1.029993303876138, 1.046501334256054, 1.01485638500529),
    @(1.02, 1.026396436944381, 1.034716328561155, 1.026693124906436, 1.043269398891151, 1, 1.033838079766376, 1.032536186168481, 1.038030357320495, 1.030034930907235, 1.046554333046655, 1.014869314952244),
    @(1.02, 1.026097914042863, 1.034450376779874, 1.026441487923372, 1.042967928011692, 1, 1.033790388290885, 1.032361442514242, 1.037829339810692, 1.029848763888608, 1.046317305083501, 1.014811478786101),
    @(1.02, 1.025127798432221, 1.033586537585917, 1.025624317266628, 1.041988632395428, 1, 1.033632727807634, 1.031792850744703, 1.037175680876073, 1.029243502424552, 1.04554662485921, 1.014623267045305),
    @(1.02, 1.024518170857396, 1.033044026508049, 1.025111247964539, 1.041373544750204, 1, 1.033531603461001, 1.031434990508731, 1.03676460583157, 1.028862946219877, 1.045062014921027, 1.014504794535626),
    @(1.02, 1.024841314909172, 1.033331562651109, 1.025383167044256, 1.041699553609821, 1, 1.033585400398631, 1.031624733082679, 1.036982532914444, 1.029064686232149, 1.045318920285263, 1.014567611830867),
    @(1.02, 1.026114183522767, 1.034464869522249, 1.026455199947275, 1.042984356652221, 1, 1.033792997575523, 1.032370968768225, 1.037840296798601, 1.029858911007066, 1.046330224629593, 1.014814631837884),
    @(1.02, 1.02759263306032, 1.035782645604486, 1.027702304119397, 1.044478000063266, 1, 1.034025178949062, 1.033235306930867, 1.038835229429024, 1.030780501530517, 1.047503504299673, 1.015100676576928)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $c = $cols[$j]
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
